$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 9 duplicates row 7 ("a1") - "980" usage_example推薦功能initial value case
$ws.Range("A7:CV7").Copy()
$ws.Range("A9").PasteSpecial(-4163)
$ws.Range("A7").Copy()
$ws.Range("A9").PasteSpecial(-4122)

# New row 10 duplicates row 8 ("b2")
$ws.Range("A8:CV8").Copy()
$ws.Range("A10").PasteSpecial(-4163)
$ws.Range("A8").Copy()
$ws.Range("A10").PasteSpecial(-4122)

$excel.CutCopyMode = $false
